$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1250612
$ws.Range("J9").Value = 4999999.5
$ws.Range("L9").Value = 4999999.5
$ws.Range("N9").Value = -5000337.5

$ws.Range("H40").Value = 5542.9287
$ws.Range("I40").Value = 2800
$ws.Range("J40").Value = 6640.1
$ws.Range("K40").Value = 2800
$ws.Range("L40").Value = 6640.1
$ws.Range("M40").Value = -2625
$ws.Range("N40").Value = -6990.1

$ws.Range("H53").Value = 390.26086
$ws.Range("I53").Value = 450.29413
$ws.Range("J53").Value = 220.16667
$ws.Range("K53").Value = 450.29413
$ws.Range("L53").Value = 220.16667
$ws.Range("M53").Value = 186.70587
$ws.Range("N53").Value = -1494.16667

$ws.Range("H62").Value = 2728.75
$ws.Range("I62").Value = 2728.75
$ws.Range("K62").Value = 2728.75
$ws.Range("M62").Value = -2104.75

$ws.Range("H65").Value = 2728.75
$ws.Range("I65").Value = 2728.75
$ws.Range("K65").Value = 13643.75
$ws.Range("M65").Value = -10523.75

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H132").Value = 36839.035
$ws.Range("I132").Value = 2362.913
$ws.Range("K132").Value = 7088.739
$ws.Range("M132").Value = -4558.739

$ws.Range("H138").Value = 5795.359
$ws.Range("I138").Value = 1497
$ws.Range("J138").Value = 5908.4736
$ws.Range("K138").Value = 4491
$ws.Range("L138").Value = 17725.4208
$ws.Range("M138").Value = 649
$ws.Range("N138").Value = -28005.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1143.6111
$ws.Range("I45").Value = 653.2
$ws.Range("K45").Value = 653.2
$ws.Range("M45").Value = -276.2

$ws.Range("H74").Value = 2279.875
$ws.Range("I74").Value = 2291.2856
$ws.Range("J74").Value = 2200
$ws.Range("K74").Value = 2291.2856
$ws.Range("L74").Value = 2200
$ws.Range("M74").Value = -1417.2856
$ws.Range("N74").Value = -3948

$ws.Range("H77").Value = 2279.875
$ws.Range("I77").Value = 2291.2856
$ws.Range("J77").Value = 2200
$ws.Range("K77").Value = 11456.428
$ws.Range("L77").Value = 11000
$ws.Range("M77").Value = -7088.428
$ws.Range("N77").Value = -19736

$ws.Range("H102").Value = 7019.6
$ws.Range("I102").Value = 7039.8
$ws.Range("J102").Value = 6999.4
$ws.Range("K102").Value = 7039.8
$ws.Range("L102").Value = 6999.4
$ws.Range("M102").Value = -5417.8
$ws.Range("N102").Value = -10243.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5217.143
$ws.Range("I22").Value = 170
$ws.Range("J22").Value = 6058.3335
$ws.Range("K22").Value = 170
$ws.Range("L22").Value = 6058.3335
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = -6404.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 772
$ws.Range("I7").Value = 135.2
$ws.Range("K7").Value = 135.2
$ws.Range("M7").Value = -22.19999999999999

$ws.Range("H31").Value = 4432.353
$ws.Range("I31").Value = 1538.3334
$ws.Range("J31").Value = 5052.5
$ws.Range("K31").Value = 1538.3334
$ws.Range("L31").Value = 5052.5
$ws.Range("M31").Value = -1243.3334
$ws.Range("N31").Value = -5642.5

$ws.Range("H34").Value = 4432.353
$ws.Range("I34").Value = 1538.3334
$ws.Range("J34").Value = 5052.5
$ws.Range("K34").Value = 1538.3334
$ws.Range("L34").Value = 5052.5
$ws.Range("M34").Value = -1336.3334
$ws.Range("N34").Value = -5456.5

$ws.Range("H62").Value = 10958.385
$ws.Range("I62").Value = 11704.917
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 11704.917
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -11080.917
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 10958.385
$ws.Range("I65").Value = 11704.917
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 58524.585
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -55404.585
$ws.Range("N65").Value = -16240

$ws.Range("H105").Value = 3908.7144
$ws.Range("I105").Value = 1750
$ws.Range("K105").Value = 1750
$ws.Range("M105").Value = -3

$ws.Range("H107").Value = 1040.1111
$ws.Range("I107").Value = 787.4545000000001
$ws.Range("K107").Value = 787.4545000000001
$ws.Range("M107").Value = 1132.5455

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 4843.1113
$ws.Range("I132").Value = 3481
$ws.Range("J132").Value = 5932.8
$ws.Range("K132").Value = 10443
$ws.Range("L132").Value = 17798.4
$ws.Range("M132").Value = -7913
$ws.Range("N132").Value = -22858.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5062.4688
$ws.Range("I56").Value = 5062.4688
$ws.Range("K56").Value = 5062.4688
$ws.Range("M56").Value = -4532.4688

$ws.Range("H68").Value = 1925
$ws.Range("J68").Value = 1925
$ws.Range("L68").Value = 5775
$ws.Range("N68").Value = -7397

$ws.Range("H71").Value = 1925
$ws.Range("J71").Value = 1925
$ws.Range("L71").Value = 17325
$ws.Range("N71").Value = -25437

$ws.Range("H92").Value = 278
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H97").Value = 1035
$ws.Range("I97").Value = 1059.7273
$ws.Range("J97").Value = 899
$ws.Range("K97").Value = 3179.1819
$ws.Range("L97").Value = 2697
$ws.Range("M97").Value = -2683.1819
$ws.Range("N97").Value = -3689

$ws.Range("H98").Value = 375
$ws.Range("I98").Value = 375
$ws.Range("K98").Value = 1125
$ws.Range("M98").Value = 373

$ws.Range("H121").Value = 909.73334
$ws.Range("I121").Value = 531.1667
$ws.Range("J121").Value = 1162.1111
$ws.Range("K121").Value = 1593.5001
$ws.Range("L121").Value = 3486.3333
$ws.Range("M121").Value = -283.5001
$ws.Range("N121").Value = -6106.3333

$ws.Range("H124").Value = 3008.6667
$ws.Range("J124").Value = 2997
$ws.Range("L124").Value = 8991
$ws.Range("N124").Value = -18811

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H131").Value = 28149648
$ws.Range("J131").Value = 83333336
$ws.Range("L131").Value = 250000008
$ws.Range("N131").Value = -250010088

$ws.Range("H132").Value = 1892
$ws.Range("I132").Value = 1311.75
$ws.Range("K132").Value = 11805.75
$ws.Range("M132").Value = -9275.75

$ws.Range("H137").Value = 3486.3333
$ws.Range("I137").Value = 459
$ws.Range("K137").Value = 1377
$ws.Range("M137").Value = 3723

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 366.33334
$ws.Range("I2").Value = 287.125
$ws.Range("K2").Value = 287.125
$ws.Range("M2").Value = -174.125

$ws.Range("H46").Value = 1987.5
$ws.Range("I46").Value = 1987.5
$ws.Range("K46").Value = 1987.5
$ws.Range("M46").Value = -1831.5

$ws.Range("H70").Value = 8934.362999999999
$ws.Range("J70").Value = 10397.143
$ws.Range("L70").Value = 10397.143
$ws.Range("N70").Value = -10937.143

$ws.Range("H73").Value = 8934.362999999999
$ws.Range("J73").Value = 10397.143
$ws.Range("L73").Value = 10397.143
$ws.Range("N73").Value = -12269.143

$ws.Range("H122").Value = 2429.889
$ws.Range("I122").Value = 2553.1428
$ws.Range("J122").Value = 1998.5
$ws.Range("K122").Value = 7659.428400000001
$ws.Range("L122").Value = 5995.5
$ws.Range("M122").Value = -5209.428400000001
$ws.Range("N122").Value = -10895.5

$ws.Range("H126").Value = 2819.9473
$ws.Range("I126").Value = 2542.7778
$ws.Range("J126").Value = 3069.4
$ws.Range("K126").Value = 7628.3334
$ws.Range("L126").Value = 9208.200000000001
$ws.Range("M126").Value = -5158.3334
$ws.Range("N126").Value = -14148.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2225.4443
$ws.Range("I16").Value = 2839
$ws.Range("K16").Value = 2839
$ws.Range("M16").Value = -2669

$ws.Range("H40").Value = 5202.5386
$ws.Range("J40").Value = 5707.75
$ws.Range("L40").Value = 5707.75
$ws.Range("N40").Value = -5979.75

$ws.Range("H46").Value = 2008
$ws.Range("I46").Value = 912.1667
$ws.Range("J46").Value = 2513.7693
$ws.Range("K46").Value = 912.1667
$ws.Range("L46").Value = 2513.7693
$ws.Range("M46").Value = -724.1667
$ws.Range("N46").Value = -2889.7693

$ws.Range("H55").Value = 215.5625
$ws.Range("I55").Value = 337.8889
$ws.Range("J55").Value = 58.285713
$ws.Range("K55").Value = 337.8889
$ws.Range("L55").Value = 58.285713
$ws.Range("M55").Value = -164.8889
$ws.Range("N55").Value = -404.285713

$ws.Range("H122").Value = 4519.148
$ws.Range("I122").Value = 3808.6
$ws.Range("J122").Value = 6549.2856
$ws.Range("K122").Value = 11425.8
$ws.Range("L122").Value = 19647.8568
$ws.Range("M122").Value = -8975.799999999999
$ws.Range("N122").Value = -24547.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 24000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H54").Value = 25250
$ws.Range("J54").Value = 11000
$ws.Range("L54").Value = 11000
$ws.Range("N54").Value = -12040

$ws.Range("H107").Value = 888.0833
$ws.Range("I107").Value = 888.0833
$ws.Range("K107").Value = 2664.2499
$ws.Range("M107").Value = -744.2498999999998

$ws.Range("H126").Value = 3067.3333
$ws.Range("I126").Value = 2489
$ws.Range("J126").Value = 6248.1665
$ws.Range("K126").Value = 7467
$ws.Range("L126").Value = 18744.4995
$ws.Range("M126").Value = -4997
$ws.Range("N126").Value = -23684.4995

$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
